# The "Förändrad" (Changed) column C holds a date serial number that was
# bumped by one day (45177 -> 45178) for every data row (rows 2-398) on the
# single worksheet "Avverkningsanmälningar".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 398 }

$ws.Range("C2:C$lastRow").Value = 45178
